$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    @("12÷9=1, 3", "41÷5=8, 1", "61÷3=20, 1", "71÷4=17, 3", "36÷8=4, 4"),
    @("80÷3=26, 2", "87÷2=43, 1", "19÷2=9, 1", "60÷7=8, 4", "10÷8=1, 2"),
    @("42÷2=21, 0", "86÷3=28, 2", "92÷3=30, 2", "84÷4=21, 0", "70÷7=10, 0"),
    @("60÷6=10, 0", "11÷4=2, 3", "15÷9=1, 6", "18÷2=9, 0", "98÷7=14, 0"),
    @("88÷9=9, 7", "23÷5=4, 3", "60÷2=30, 0", "57÷5=11, 2", "69÷8=8, 5")
)

$rowIndexes = @(1, 5, 9, 13, 17)

for ($i = 0; $i -lt $rowIndexes.Length; $i++) {
    $row = $rowIndexes[$i]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $newValues[$i][$col - 1]
    }
}
